$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helpers -----------------------------------------------------------
# Write a value that must stay literal TEXT (dates like "2022-09-01" and
# big numeric-looking ids would otherwise be auto-converted by Excel into
# a date serial / scientific-notation number). Pre-formatting the cell as
# Text ("@") keeps the original string intact.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
}

# --- New header column H (copy bold header style from A1, then set text)
$ws.Range("A1").Copy($ws.Range("H1")) | Out-Null
$ws.Range("H1").Value = "class_year"

# --- Row 2 : new "halfyear" period (id=4) -------------------------------
$ws.Range("A2").Value = 4
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "halfyear"
Set-TextValue "E2" "2022-09-01"
Set-TextValue "F2" "2022-12-31"
Set-TextValue "G2" "1985861362764370224"
$ws.Range("H2").Value = "10,11,12,13"

# --- Row 3 : new "halfyear" period (id=5) -------------------------------
$ws.Range("A3").Value = 5
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "halfyear"
Set-TextValue "E3" "2023-01-01"
Set-TextValue "F3" "2023-06-03"
Set-TextValue "G3" "1985861362764370225"
$ws.Range("H3").Value = "10,11,12,13"

# --- Row 4 : moved "trimester" period (id=1), now with class_year ------
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = "trimester"
Set-TextValue "E4" "2022-09-01"
Set-TextValue "F4" "2022-11-30"
Set-TextValue "G4" "1985862771513643316"
$ws.Range("H4").Value = "5,6,7,8,9"

# --- Row 5 : moved "trimester" period (id=2), now with class_year ------
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = "trimester"
Set-TextValue "E5" "2022-12-01"
Set-TextValue "F5" "2023-02-28"
Set-TextValue "G5" "1985862771513643317"
$ws.Range("H5").Value = "5,6,7,8,9"

# --- Row 6 : moved "trimester" period (id=3), now with class_year ------
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = "trimester"
Set-TextValue "E6" "2023-03-01"
Set-TextValue "F6" "2023-06-02"
Set-TextValue "G6" "1985862771513643318"
$ws.Range("H6").Value = "5,6,7,8,9"
